$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 71; this shifts the existing rows 71..200
# down by two (they reappear at 73..202), carrying their formatting with them
# (e.g. row 71's date-column style now applies to row 73, etc.).
$ws.Rows("71:72").Insert()

# Populate the two newly inserted rows (71 and 72) with the new weekly record.
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = 'Vega Monumental Concepción'
$ws.Range("C71").Value = 'Bíobío'
$ws.Range("D71").Value = 44467
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112006
$ws.Range("G71").Value = 'Repollo'
$ws.Range("H71").Value = 'Crespo record'
$ws.Range("I71").Value = 'Primera'
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 600
$ws.Range("L71").Value = 700
$ws.Range("M71").Value = 650
$ws.Range("N71").Value = '$/unidad'
$ws.Range("O71").Value = 'Región Metropolitana'
$ws.Range("P71").Value = 650
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = 'Hortaliza'

$ws.Range("A72").Value = 11
$ws.Range("B72").Value = 'Vega Monumental Concepción'
$ws.Range("C72").Value = 'Bíobío'
$ws.Range("D72").Value = 44467
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 100112006
$ws.Range("G72").Value = 'Repollo'
$ws.Range("H72").Value = 'Crespo record'
$ws.Range("I72").Value = 'Segunda'
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 500
$ws.Range("L72").Value = 500
$ws.Range("M72").Value = 500
$ws.Range("N72").Value = '$/unidad'
$ws.Range("O72").Value = 'Región Metropolitana'
$ws.Range("P72").Value = 500
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = 'Hortaliza'
